{"js": "// Replace the date line and each two-digit-by-two-digit multiplication\n// problem in the worksheet with its new value, per the commit's diff.\n// Every \"old\" text below is unique in the document, so an exact,\n// case-sensitive whole-match search safely targets only the one run\n// that should change.\nconst replacements = [\n  [\"2025-05-12 Monday\", \"2025-05-13 Tuesday\"],\n  [\"61\u00d717=\", \"64\u00d798=\"],\n  [\"93\u00d750=\", \"38\u00d767=\"],\n  [\"68\u00d771=\", \"41\u00d716=\"],\n  [\"98\u00d755=\", \"44\u00d758=\"],\n  [\"64\u00d725=\", \"63\u00d734=\"],\n  [\"35\u00d725=\", \"36\u00d742=\"],\n  [\"17\u00d718=\", \"62\u00d752=\"],\n  [\"98\u00d753=\", \"33\u00d711=\"],\n  [\"80\u00d735=\", \"40\u00d766=\"],\n  [\"88\u00d776=\", \"13\u00d754=\"],\n  [\"16\u00d792=\", \"47\u00d741=\"],\n  [\"30\u00d749=\", \"83\u00d716=\"],\n  [\"37\u00d760=\", \"31\u00d719=\"],\n  [\"92\u00d780=\", \"60\u00d789=\"],\n  [\"80\u00d783=\", \"13\u00d766=\"],\n  [\"75\u00d711=\", \"21\u00d777=\"],\n  [\"81\u00d794=\", \"33\u00d795=\"],\n  [\"85\u00d742=\", \"69\u00d763=\"],\n  [\"41\u00d769=\", \"12\u00d733=\"],\n  [\"16\u00d794=\", \"96\u00d779=\"],\n  [\"17\u00d764=\", \"53\u00d782=\"],\n  [\"87\u00d720=\", \"46\u00d782=\"],\n  [\"90\u00d721=\", \"23\u00d785=\"],\n  [\"93\u00d735=\", \"97\u00d745=\"],\n  [\"20\u00d757=\", \"71\u00d733=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each two-digit-by-two-digit multiplication\n# problem in the worksheet with its new value, per the commit's diff.\n# Every \"old\" text below is unique in the document, so an exact,\n# case-sensitive Find/Replace (wdReplaceAll, bounded to a single hit)\n# safely targets only the one run that should change.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-05-12 Monday\", \"2025-05-13 Tuesday\"),\n    @(\"61\u00d717=\", \"64\u00d798=\"),\n    @(\"93\u00d750=\", \"38\u00d767=\"),\n    @(\"68\u00d771=\", \"41\u00d716=\"),\n    @(\"98\u00d755=\", \"44\u00d758=\"),\n    @(\"64\u00d725=\", \"63\u00d734=\"),\n    @(\"35\u00d725=\", \"36\u00d742=\"),\n    @(\"17\u00d718=\", \"62\u00d752=\"),\n    @(\"98\u00d753=\", \"33\u00d711=\"),\n    @(\"80\u00d735=\", \"40\u00d766=\"),\n    @(\"88\u00d776=\", \"13\u00d754=\"),\n    @(\"16\u00d792=\", \"47\u00d741=\"),\n    @(\"30\u00d749=\", \"83\u00d716=\"),\n    @(\"37\u00d760=\", \"31\u00d719=\"),\n    @(\"92\u00d780=\", \"60\u00d789=\"),\n    @(\"80\u00d783=\", \"13\u00d766=\"),\n    @(\"75\u00d711=\", \"21\u00d777=\"),\n    @(\"81\u00d794=\", \"33\u00d795=\"),\n    @(\"85\u00d742=\", \"69\u00d763=\"),\n    @(\"41\u00d769=\", \"12\u00d733=\"),\n    @(\"16\u00d794=\", \"96\u00d779=\"),\n    @(\"17\u00d764=\", \"53\u00d782=\"),\n    @(\"87\u00d720=\", \"46\u00d782=\"),\n    @(\"90\u00d721=\", \"23\u00d785=\"),\n    @(\"93\u00d735=\", \"97\u00d745=\"),\n    @(\"20\u00d757=\", \"71\u00d733=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Execute(\n        [ref]$oldText,\n        [ref]$true,\n        [ref]$false,\n        [ref]$false,\n        [ref]$false,\n        [ref]$false,\n        [ref]$true,\n        [ref]1,\n        [ref]$false,\n        [ref]$newText,\n        [ref]2\n    )\n}\n"}
